$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the new rows of working-hours data (dates + Developer1 hours)
# Copy the date formatting (style "s=1", numFmtId 14) from an existing date
# cell so no new style gets minted, then overwrite the value.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 44304
$ws.Range("B11").Value = "1"

$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 44307
$ws.Range("B12").Value = "2.5"

$ws.Range("A10").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 44308
$ws.Range("B13").Value = "4"

$ws.Range("A10").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 44309
$ws.Range("B14").Value = "4"

# Extend the sheet with one more blank data row, matching the style of the
# row above it.
$ws.Range("B29:C29").Copy()
$ws.Range("B30:C30").PasteSpecial(-4122)

# Update selection to reflect where the user was working
$ws.Range("I21").Select()
